$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settlement")

# Correct the settlement denomination counts and recompute the cash value
$ws.Range("B2").Value = "2"
$ws.Range("B3").Value = "8"
$ws.Range("B8").Value = "7842.08"
$ws.Range("B4").Value = "11"
$ws.Range("B6").Value = "19"

# "Settlement" (not "Load In") is the correct active sheet/tab
$ws.Activate()
